$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9127805829048157
$ws.Range("B1").Value = 1.636133670806885
$ws.Range("C1").Value = 2.166375637054443
$ws.Range("D1").Value = 2.25888729095459
$ws.Range("E1").Value = 1.557436108589172
